$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Row 2
$r = $ws.Range("D2")
$r.Value = "'60.271.14"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +2.20%  "

# Row 3
$r = $ws.Range("D3")
$r.Value = "'2.585.83"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +1.88%  "

# Row 4
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$r = $ws.Range("D5")
$r.Value = "'505.40"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "

# Row 6
$r = $ws.Range("D6")
$r.Value = "'152.42"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -3.75%  "

# Row 7
$ws.Range("E7").Value = "  +0.48%  "

# Row 8
$r = $ws.Range("D8")
$r.Value = "'0.579"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -7.00%  "

# Row 9
$r = $ws.Range("D9")
$r.Value = "'2.591.03"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +0.66%  "

# Row 10
$r = $ws.Range("D10")
$r.Value = "'6.68"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +7.90%  "

# Row 11
$ws.Range("E11").Value = "  -0.81%  "

# Row 12
$ws.Range("E12").Value = "  +0.62%  "

# Row 13
$ws.Range("E13").Value = "  +0.84%  "

# Row 14
$r = $ws.Range("D14")
$r.Value = "'3.039.02"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +1.79%  "

# Row 15
$r = $ws.Range("D15")
$r.Value = "'60.262.80"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +2.36%  "

# Row 16
$r = $ws.Range("D16")
$r.Value = "'21.53"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -2.07%  "

# Row 17
$ws.Range("E17").Value = "  +0.99%  "

# Row 18
$r = $ws.Range("D18")
$r.Value = "'2.591.30"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19
$ws.Range("E19").Value = "  +1.13%  "

# Row 20
$r = $ws.Range("D20")
$r.Value = "'345.76"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +3.05%  "

# Row 21
$r = $ws.Range("D21")
$r.Value = "'10.29"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -0.72%  "

# Row 22
$ws.Range("E22").Value = "  +0.50%  "

# Row 23
$r = $ws.Range("D23")
$r.Value = "'0.997"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "

# Row 24
$r = $ws.Range("D24")
$r.Value = "'60.26"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "

# Row 25
$ws.Range("E25").Value = "  +0.42%  "

# Row 26
$ws.Range("E26").Value = "  -1.42%  "

# Row 27
$r = $ws.Range("D27")
$r.Value = "'2.697.79"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +1.57%  "

# Row 28
$r = $ws.Range("D28")
$r.Value = "'0.997"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +0.23%  "

# Row 29
$r = $ws.Range("D29")
$r.Value = "'0.0₃0839"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +1.11%  "

# Row 30
$r = $ws.Range("D30")
$r.Value = "'7.45"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.66%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$r = $ws.Range("D32")
$r.Value = "'19.28"
$r.Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "

# Row 33
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$r = $ws.Range("D33")
$r.Value = "'154.49"
$r.Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "

# Row 34
$ws.Range("E34").Value = "  -1.14%  "

# Row 35
$r = $ws.Range("D35")
$r.Value = "'5.73"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +4.35%  "

# Row 36
$r = $ws.Range("D36")
$r.Value = "'3.98"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +1.71%  "

# Row 37
$ws.Range("E37").Value = "  -0.63%  "

# Row 38
$r = $ws.Range("D38")
$r.Value = "'0.857"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +19.37%  "

# Row 39
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$r = $ws.Range("D39")
$r.Value = "'0.845"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -1.32%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$r = $ws.Range("D40")
$r.Value = "'3.77"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$r = $ws.Range("D41")
$r.Value = "'35.79"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +2.74%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$r = $ws.Range("D42")
$r.Value = "'297.73"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +2.45%  "

# Row 43
$r = $ws.Range("D43")
$r.Value = "'1.44"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +0.81%  "

# Row 44
$r = $ws.Range("D44")
$r.Value = "'0.618"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -1.20%  "

# Row 45
$r = $ws.Range("D45")
$r.Value = "'0.0997"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -1.97%  "

# Row 46
$ws.Range("E46").Value = "  -0.29%  "

# Row 47
$r = $ws.Range("D47")
$r.Value = "'0.998"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "

# Row 48
$r = $ws.Range("D48")
$r.Value = "'19.65"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "

# Row 49
$r = $ws.Range("D49")
$r.Value = "'4.86"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "

# Row 50
$ws.Range("E50").Value = "  -2.34%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$r = $ws.Range("D51")
$r.Value = "'10.30"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
